$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pthlh"
$ws.Range("C2").Value = "Pth2r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2032693333333333
$ws.Range("H2").Value = 0.609808
$ws.Range("I2").Value = 0.01014574830092973
$ws.Range("J2").Value = 0.01014574830092973
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01042033333333333
$ws.Range("N2").Value = 0.031261
$ws.Range("O2").Value = 0.02238908055155877
$ws.Range("P2").Value = 0.02238908055155877
$ws.Range("Q2").Value = 0.002118134209777777
$ws.Range("R2").Value = 0.019063207888
$ws.Range("S2").Value = 0.0002271539759653562
$ws.Range("T2").Value = 0.0002271539759653562

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pthlh"
$ws.Range("C3").Value = "Pth2r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2032693333333333
$ws.Range("H3").Value = 0.609808
$ws.Range("I3").Value = 0.01014574830092973
$ws.Range("J3").Value = 0.01014574830092973
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.09745433333333332
$ws.Range("N3").Value = 0.292363
$ws.Range("O3").Value = 0.2093899349763403
$ws.Range("P3").Value = 0.2093899349763404
$ws.Range("Q3").Value = 0.01980947736711111
$ws.Range("R3").Value = 0.178285296304
$ws.Range("S3").Value = 0.002124417577017992
$ws.Range("T3").Value = 0.002124417577017992

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pthlh"
$ws.Range("C4").Value = "Pth2r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2032693333333333
$ws.Range("H4").Value = 0.609808
$ws.Range("I4").Value = 0.01014574830092973
$ws.Range("J4").Value = 0.01014574830092973
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3575456666666667
$ws.Range("N4").Value = 1.072637
$ws.Range("O4").Value = 0.7682209844721009
$ws.Range("P4").Value = 0.7682209844721009
$ws.Range("Q4").Value = 0.07267806929955556
$ws.Range("R4").Value = 0.654102623696
$ws.Range("S4").Value = 0.007794176747946381
$ws.Range("T4").Value = 0.007794176747946381

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Pthlh"
$ws.Range("C5").Value = "Pth2r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.94203033333333
$ws.Range("H5").Value = 56.826091
$ws.Range("I5").Value = 0.9454503978493691
$ws.Range("J5").Value = 0.9454503978493691
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01042033333333333
$ws.Range("N5").Value = 0.031261
$ws.Range("O5").Value = 0.02238908055155877
$ws.Range("P5").Value = 0.02238908055155877
$ws.Range("Q5").Value = 0.1973822700834444
$ws.Range("R5").Value = 1.776440430751
$ws.Range("S5").Value = 0.02116776511495281
$ws.Range("T5").Value = 0.02116776511495281

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pthlh"
$ws.Range("C6").Value = "Pth2r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 18.94203033333333
$ws.Range("H6").Value = 56.826091
$ws.Range("I6").Value = 0.9454503978493691
$ws.Range("J6").Value = 0.9454503978493691
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.09745433333333332
$ws.Range("N6").Value = 0.292363
$ws.Range("O6").Value = 0.2093899349763403
$ws.Range("P6").Value = 0.2093899349763404
$ws.Range("Q6").Value = 1.845982938114777
$ws.Range("R6").Value = 16.613846443033
$ws.Range("S6").Value = 0.1979677973290345
$ws.Range("T6").Value = 0.1979677973290345

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pthlh"
$ws.Range("C7").Value = "Pth2r"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 18.94203033333333
$ws.Range("H7").Value = 56.826091
$ws.Range("I7").Value = 0.9454503978493691
$ws.Range("J7").Value = 0.9454503978493691
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3575456666666667
$ws.Range("N7").Value = 1.072637
$ws.Range("O7").Value = 0.7682209844721009
$ws.Range("P7").Value = 0.7682209844721009
$ws.Range("Q7").Value = 6.772640863551889
$ws.Range("R7").Value = 60.953767771967
$ws.Range("S7").Value = 0.7263148354053818
$ws.Range("T7").Value = 0.7263148354053818

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Pthlh"
$ws.Range("C8").Value = "Pth2r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.889628
$ws.Range("H8").Value = 2.668884
$ws.Range("I8").Value = 0.04440385384970112
$ws.Range("J8").Value = 0.04440385384970111
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01042033333333333
$ws.Range("N8").Value = 0.031261
$ws.Range("O8").Value = 0.02238908055155877
$ws.Range("P8").Value = 0.02238908055155877
$ws.Range("Q8").Value = 0.009270220302666665
$ws.Range("R8").Value = 0.08343198272399999
$ws.Range("S8").Value = 0.0009941614606406011
$ws.Range("T8").Value = 0.0009941614606406011

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Pthlh"
$ws.Range("C9").Value = "Pth2r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.889628
$ws.Range("H9").Value = 2.668884
$ws.Range("I9").Value = 0.04440385384970112
$ws.Range("J9").Value = 0.04440385384970111
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.09745433333333332
$ws.Range("N9").Value = 0.292363
$ws.Range("O9").Value = 0.2093899349763403
$ws.Range("P9").Value = 0.2093899349763404
$ws.Range("Q9").Value = 0.08669810365466665
$ws.Range("R9").Value = 0.7802829328919999
$ws.Range("S9").Value = 0.009297720070287837
$ws.Range("T9").Value = 0.009297720070287835

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Pthlh"
$ws.Range("C10").Value = "Pth2r"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.889628
$ws.Range("H10").Value = 2.668884
$ws.Range("I10").Value = 0.04440385384970112
$ws.Range("J10").Value = 0.04440385384970111
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3575456666666667
$ws.Range("N10").Value = 1.072637
$ws.Range("O10").Value = 0.7682209844721009
$ws.Range("P10").Value = 0.7682209844721009
$ws.Range("Q10").Value = 0.3180826363453333
$ws.Range("R10").Value = 2.862743727108
$ws.Range("S10").Value = 0.03411197231877268
$ws.Range("T10").Value = 0.03411197231877267
